$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Row 2
Set-TextCell $ws.Range("D2") "30.117.11"
$ws.Range("E2").Value = "  -1.20%  "

# Row 3
Set-TextCell $ws.Range("D3") "1.856.28"
$ws.Range("E3").Value = "  -2.89%  "

# Row 4
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
Set-TextCell $ws.Range("D5") "233.92"
$ws.Range("E5").Value = "  -2.29%  "

# Row 6
$ws.Range("E6").Value = "  +0.09%  "

# Row 7
Set-TextCell $ws.Range("D7") "0.4693"
$ws.Range("E7").Value = "  -1.72%  "

# Row 8
Set-TextCell $ws.Range("D8") "0.2816"
$ws.Range("E8").Value = "  -1.05%  "

# Row 9
Set-TextCell $ws.Range("D9") "0.06542"
$ws.Range("E9").Value = "  -2.49%  "

# Row 10
Set-TextCell $ws.Range("D10") "20.13"
$ws.Range("E10").Value = "  +3.11%  "

# Row 11
Set-TextCell $ws.Range("D11") "0.07797"
$ws.Range("E11").Value = "  +0.54%  "

# Row 12
Set-TextCell $ws.Range("D12") "97.15"
$ws.Range("E12").Value = "  -6.00%  "

# Row 13
Set-TextCell $ws.Range("D13") "1.854.51"
$ws.Range("E13").Value = "  -2.98%  "

# Row 14
Set-TextCell $ws.Range("D14") "5.080"
$ws.Range("E14").Value = "  -2.09%  "

# Row 15
Set-TextCell $ws.Range("D15") "0.6702"
$ws.Range("E15").Value = "  +0.06%  "

# Row 16
Set-TextCell $ws.Range("D16") "284.75"
$ws.Range("E16").Value = "  +3.15%  "

# Row 17
Set-TextCell $ws.Range("D17") "30.154.48"
$ws.Range("E17").Value = "  -1.11%  "

# Row 18
Set-TextCell $ws.Range("D18") "0.9996"
$ws.Range("E18").Value = "  +0.05%  "

# Row 19
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
Set-TextCell $ws.Range("D19") "5.422"
$ws.Range("E19").Value = "  +0.80%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
Set-TextCell $ws.Range("D20") "12.60"
$ws.Range("E20").Value = "  -0.05%  "

# Row 21
Set-TextCell $ws.Range("D21") "2.104.14"
$ws.Range("E21").Value = "  -2.43%  "

# Row 22
Set-TextCell $ws.Range("D22") "0.000007231"
$ws.Range("E22").Value = "  -3.46%  "

# Row 23
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
Set-TextCell $ws.Range("D24") "6.145"
$ws.Range("E24").Value = "  -2.41%  "

# Row 25
Set-TextCell $ws.Range("D25") "167.78"

# Row 26
Set-TextCell $ws.Range("D26") "9.298"
$ws.Range("E26").Value = "  -0.56%  "

# Row 27
Set-TextCell $ws.Range("D27") "19.05"
$ws.Range("E27").Value = "  -0.86%  "

# Row 28
Set-TextCell $ws.Range("D28") "1.926"
$ws.Range("E28").Value = "  -7.36%  "

# Row 29
Set-TextCell $ws.Range("D29") "1.339"
$ws.Range("E29").Value = "  -3.12%  "

# Row 30
Set-TextCell $ws.Range("D30") "0.09648"
$ws.Range("E30").Value = "  -3.37%  "

# Row 31
Set-TextCell $ws.Range("D31") "4.401"
$ws.Range("E31").Value = "  -4.16%  "

# Row 32
Set-TextCell $ws.Range("D32") "1.469"
$ws.Range("E32").Value = "  -2.80%  "

# Row 33
Set-TextCell $ws.Range("D33") "4.097"
$ws.Range("E33").Value = "  -3.70%  "

# Row 34
Set-TextCell $ws.Range("D34") "0.04668"
$ws.Range("E34").Value = "  -0.86%  "

# Row 35
Set-TextCell $ws.Range("D35") "0.6968"
$ws.Range("E35").Value = "  -4.27%  "

# Row 36
Set-TextCell $ws.Range("D36") "1.088"
$ws.Range("E36").Value = "  -2.36%  "

# Row 37
Set-TextCell $ws.Range("D37") "0.9995"
$ws.Range("E37").Value = "  +0.13%  "

# Row 38
Set-TextCell $ws.Range("D38") "2.701"
$ws.Range("E38").Value = "  -0.29%  "

# Row 39
Set-TextCell $ws.Range("D39") "0.01859"
$ws.Range("E39").Value = "  -2.46%  "

# Row 40
Set-TextCell $ws.Range("D40") "6.311"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41
Set-TextCell $ws.Range("D41") "2.506"
$ws.Range("E41").Value = "  -3.90%  "

# Row 42
Set-TextCell $ws.Range("D42") "71.81"
$ws.Range("E42").Value = "  -3.44%  "

# Row 43
Set-TextCell $ws.Range("D43") "0.8610"
$ws.Range("E43").Value = "  -0.05%  "

# Row 44
Set-TextCell $ws.Range("D44") "1.938"
$ws.Range("E44").Value = "  -0.97%  "

# Row 45
Set-TextCell $ws.Range("D45") "104.29"
$ws.Range("E45").Value = "  -2.14%  "

# Row 46
$ws.Range("B46").Value = "TheSandbox"
$ws.Range("C46").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextCell $ws.Range("D46") "0.4164"
$ws.Range("E46").Value = "  -2.27%  "

# Row 47
$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
Set-TextCell $ws.Range("D47") "0.9998"
$ws.Range("E47").Value = "  +0.06%  "

# Row 48
Set-TextCell $ws.Range("D48") "1.019.27"
$ws.Range("E48").Value = "  +7.17%  "

# Row 49
Set-TextCell $ws.Range("D49") "7.246"
$ws.Range("E49").Value = "  -2.19%  "

# Row 50
Set-TextCell $ws.Range("D50") "9.174"
$ws.Range("E50").Value = "  +4.46%  "

# Row 51
Set-TextCell $ws.Range("D51") "33.74"
$ws.Range("E51").Value = "  -2.49%  "
